# Insert three "section break" paragraphs into the body, reproducing the
# target OOXML diff:
#   1. after "...Violin plot..." paragraph  -> new sectPr, default (next page) type
#   2. after "...Temporal trends..." paragraph -> new sectPr, type="continuous"
#   3. after "...Joint distribution..." paragraph -> new sectPr, type="continuous"
# and the trailing body-level <w:sectPr> gains <w:type w:val="continuous"/>.
#
# wdSectionBreakNextPage = 2, wdSectionBreakContinuous = 3 (Range.InsertBreak)
# wdSectionContinuous = 0, wdSectionNewPage = 2 (PageSetup.SectionStart)

$d = $word.ActiveDocument

function Find-ParagraphStartingWith($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "$needle*") {
            return $i
        }
    }
    return -1
}

# --- Step 1: break after the "Violin plot..." paragraph (next-page section) ---
$idx = Find-ParagraphStartingWith "Violin plot showing value distributions"
$d.Paragraphs.Item($idx).Range.InsertBreak(2)

# --- Step 2: break after the "Temporal trends..." paragraph (continuous) ---
$idx = Find-ParagraphStartingWith "Temporal trends across different groups"
$d.Paragraphs.Item($idx).Range.InsertBreak(3)

# --- Step 3: break after the "Joint distribution..." paragraph (continuous) ---
$idx = Find-ParagraphStartingWith "Joint distribution of Value1 and Value2"
$d.Paragraphs.Item($idx).Range.InsertBreak(3)

# Inserting the plain next-page break in step 1 resets the "start type" of the
# section immediately following it back to next-page, clobbering the
# continuous type that step 2 had set. Put it back to continuous explicitly
# so the section that begins right after "Violin plot..." picture / ends at
# the "Temporal trends..." break is marked continuous, matching the target.
for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $s = $d.Sections.Item($i)
    $txt = $s.Range.Text
    if ($txt -like "*Time Series Analysis*") {
        $s.PageSetup.SectionStart = 0
    }
}

Write-Output "Paragraphs: $($d.Paragraphs.Count); Sections: $($d.Sections.Count)"
